# Daily attendance processing - normalize "Recorded By" (column G) entries
# so that automated/system contributors are listed before human recorders.
#
# For every data row, if the "Recorded By" cell contains a comma-separated
# list of recorders that does NOT already start with "System", the list
# order is reversed (e.g. "dnasr281@gmail.com, System" becomes
# "System, dnasr281@gmail.com"). Rows already starting with "System", single
# -value rows, and empty rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    if ($val.StartsWith("System")) { continue }
    if ($val.Contains(",") -eq $false) { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $reversed = $parts[($parts.Count - 1)..0]
    $cell.Value = $reversed -join ", "
}
